$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.41859767129913905,
    -0.0099999995773742967,
    -0.0089999995688749834,
    -0.011999999884139356,
    -0.0059999995735227074,
    -0.0059999995592754374,
    -0.01999999947931741,
    -0.019999999474427099,
    -0.0059999995472459489,
    -0.021640351272147029,
    -0.0044999995513315127,
    -0.0059999995423325458,
    -0.0059999995403954287,
    -0.011999999506512538,
    0.024706585654678648,
    -0.0059999995388393401,
    -0.0059999995369652837,
    -0.0089999995197080906,
    -0.0089999995817171552,
    -0.0089999995674219235,
    -0.076844990078834385,
    -0.008999999563889638,
    -0.020254460833111665,
    -0.041999999373325103,
    -0.041999999369862095,
    -0.0059999995581669907,
    -0.0059999995567370235,
    -0.0059999995512720616,
    -0.011999999514037185,
    -0.019999999467686269,
    -0.01499999949428954,
    -0.020999999460853402,
    0.022235973211565785
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
